$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.682.87'
$ws.Cells.Item(2, 5).Value = '  +0.57%  '
$ws.Cells.Item(3, 4).Value = '3.387.16'
$ws.Cells.Item(3, 5).Value = '  -0.12%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '578.58'
$ws.Cells.Item(5, 5).Value = '  +0.65%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '137.16'
$ws.Cells.Item(6, 5).Value = '  +0.60%  '
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$ws.Cells.Item(8, 4).Value = '3.384.24'
$ws.Cells.Item(8, 5).Value = '  -0.15%  '
$ws.Cells.Item(9, 5).Value = '  -0.74%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '7.46'
$ws.Cells.Item(10, 5).Value = '  -1.69%  '
$ws.Cells.Item(11, 5).Value = '  +1.66%  '
$ws.Cells.Item(12, 5).Value = '  +0.95%  '
$ws.Cells.Item(13, 4).Value = '3.966.83'
$ws.Cells.Item(13, 5).Value = '  -0.10%  '
$ws.Cells.Item(14, 5).Value = '  +1.92%  '
$ws.Cells.Item(15, 5).Value = '  +0.91%  '
$ws.Cells.Item(16, 4).Value = '3.391.23'
$ws.Cells.Item(16, 5).Value = '  +0.08%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '25.46'
$ws.Cells.Item(17, 5).Value = '  +1.03%  '
$ws.Cells.Item(18, 4).Value = '61.831.76'
$ws.Cells.Item(18, 5).Value = '  +0.55%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '14.18'
$ws.Cells.Item(19, 5).Value = '  +0.64%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '9.51'
$ws.Cells.Item(20, 5).Value = '  +0.36%  '
$ws.Cells.Item(21, 5).Value = '  +0.33%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '378.92'
$ws.Cells.Item(22, 5).Value = '  +0.56%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.561'
$ws.Cells.Item(23, 5).Value = '  -1.48%  '
$ws.Cells.Item(24, 4).Value = '3.528.59'
$ws.Cells.Item(24, 5).Value = '  +0.11%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.998'
$ws.Cells.Item(25, 5).Value = '  -0.24%  '
$ws.Cells.Item(26, 5).Value = '  +5.14%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '71.26'
$ws.Cells.Item(27, 5).Value = '  +0.68%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.74'
$ws.Cells.Item(28, 5).Value = '  +2.47%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.62'
$ws.Cells.Item(29, 5).Value = '  -2.20%  '
$ws.Cells.Item(30, 5).Value = '  +0.09%  '
$ws.Cells.Item(31, 5).Value = '  +0.24%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.159'
$ws.Cells.Item(32, 5).Value = '  +2.11%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.17'
$ws.Cells.Item(33, 5).Value = '  +0.43%  '
$ws.Cells.Item(34, 5).Value = '  +0.06%  '
$ws.Cells.Item(35, 5).Value = '  -0.18%  '
$ws.Cells.Item(36, 2).Value = 'NEARProtocol'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.36'
$ws.Cells.Item(36, 5).Value = '  -4.33%  '
$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.56'
$ws.Cells.Item(37, 5).Value = '  -1.17%  '
$ws.Cells.Item(38, 2).Value = 'Aptos'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '6.85'
$ws.Cells.Item(38, 5).Value = '  -1.57%  '
$ws.Cells.Item(39, 2).Value = 'Monero'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '165.19'
$ws.Cells.Item(39, 5).Value = '  +1.34%  '
$ws.Cells.Item(40, 2).Value = 'Hedera'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0786'
$ws.Cells.Item(40, 5).Value = '  -0.40%  '
$ws.Cells.Item(41, 2).Value = 'Mantle'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.783'
$ws.Cells.Item(41, 5).Value = '  +2.62%  '
$ws.Cells.Item(42, 2).Value = 'ONDO'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.24'
$ws.Cells.Item(42, 5).Value = '  +1.17%  '
$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.00'
$ws.Cells.Item(43, 5).Value = '  -0.01%  '
$ws.Cells.Item(44, 2).Value = 'Stacks'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.72'
$ws.Cells.Item(44, 5).Value = '  +6.72%  '
$ws.Cells.Item(45, 2).Value = 'Filecoin'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '4.42'
$ws.Cells.Item(45, 5).Value = '  -0.06%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '24.86'
$ws.Cells.Item(46, 5).Value = '  +5.51%  '
$ws.Cells.Item(47, 2).Value = 'OKB'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '41.25'
$ws.Cells.Item(47, 5).Value = '  -0.56%  '
$ws.Cells.Item(48, 2).Value = 'Cosmos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '6.86'
$ws.Cells.Item(48, 5).Value = '  -1.60%  '
$ws.Cells.Item(49, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '22.86'
$ws.Cells.Item(49, 5).Value = '  -1.08%  '
$ws.Cells.Item(50, 2).Value = 'Maker'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(50, 4).Value = '2.341.87'
$ws.Cells.Item(50, 5).Value = '  +5.97%  '
$ws.Cells.Item(51, 2).Value = 'VeChain'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0262'
$ws.Cells.Item(51, 5).Value = '  +1.78%  '
